$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in / clear a few scattered values in the upper block (rows 19-23) ---
# Row 19 (RM 125): C19 was blank, now has a value
$ws.Range("C19").Value = 13.2

# Row 21 (RM 135): C21 had a value, now cleared
$ws.Range("C21").ClearContents()

# Row 23 (RM 140): C23 was blank, now has a value
$ws.Range("C23").Value = 12.2

# --- Remove the "RM 232" row (row 26) and the "SC 92" row (now row 27 after the
# first deletion) entirely, shifting all following rows up ---
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- After the two deletions, "SC 101" now sits at row 27 and "SC 232" at row 33 ---
# SC 101 row: B27 gains a value, C27 loses its value
$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()

# SC 232 row: C33 gains a value
$ws.Range("C33").Value = 10.4
